$d = $word.ActiveDocument

# The first paragraph currently reads "Learning github pull" split across
# three runs (with spell-check proofErr markers around "github"). Collapse
# all of that down to a single run reading "Hello githhub" using
# Find/Replace, which rewrites the whole matched range as one run and
# drops the now-stale proofErr markers.
$d.Content.Find.Execute("Learning github pull", $true, $false, $false, `
    $false, $false, $true, 1, $false, "Hello githhub", 2)

# The document also had a second, empty paragraph right after that one;
# remove it entirely so only the "Hello githhub" paragraph remains.
# Deleting both paragraph marks (the end of paragraph 1 and the end of
# paragraph 2) in one go merges the two paragraphs while keeping
# paragraph 1's own paragraph properties/identity (deleting only
# paragraph 1's mark would instead make the merged paragraph inherit
# paragraph 2's formatting).
$p1 = $d.Paragraphs.Item(1)
$p2 = $d.Paragraphs.Item(2)
$marks = $d.Range($p1.Range.End - 1, $p2.Range.End)
$marks.Delete()
